$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5332.6665
$ws.Range("I62").Value = 4924
$ws.Range("J62").Value = 6150
$ws.Range("K62").Value = 4924
$ws.Range("L62").Value = 6150
$ws.Range("M62").Value = -4300
$ws.Range("N62").Value = -7398

$ws.Range("H65").Value = 5332.6665
$ws.Range("I65").Value = 4924
$ws.Range("J65").Value = 6150
$ws.Range("K65").Value = 24620
$ws.Range("L65").Value = 30750
$ws.Range("M65").Value = -21500
$ws.Range("N65").Value = -36990

$ws.Range("H111").Value = 1099.8
$ws.Range("I111").Value = 800
$ws.Range("J111").Value = 1299.6666
$ws.Range("K111").Value = 2400
$ws.Range("L111").Value = 3898.9998
$ws.Range("M111").Value = 667
$ws.Range("N111").Value = -10032.9998

$ws.Range("H132").Value = 1435.4286
$ws.Range("I132").Value = 1302.7646
$ws.Range("J132").Value = 1999.25
$ws.Range("K132").Value = 3908.2938
$ws.Range("L132").Value = 5997.75
$ws.Range("M132").Value = -1378.2938
$ws.Range("N132").Value = -11057.75

$ws.Range("H138").Value = 1552.9166
$ws.Range("I138").Value = 1552.9166
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 4658.7498
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = 481.2502000000004
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12223.419
$ws.Range("I32").Value = 4962
$ws.Range("J32").Value = 19968.934
$ws.Range("K32").Value = 4962
$ws.Range("L32").Value = 19968.934
$ws.Range("M32").Value = -4675
$ws.Range("N32").Value = -20542.934

$ws.Range("H44").Value = 71033
$ws.Range("I44").Value = 56022
$ws.Range("J44").Value = 86044
$ws.Range("K44").Value = 56022
$ws.Range("L44").Value = 86044
$ws.Range("M44").Value = -55534
$ws.Range("N44").Value = -87020

$ws.Range("H45").Value = 10590273
$ws.Range("I45").Value = 3880
$ws.Range("J45").Value = 21176666
$ws.Range("K45").Value = 3880
$ws.Range("L45").Value = 21176666
$ws.Range("M45").Value = -3503
$ws.Range("N45").Value = -21177420

$ws.Range("H61").Value = 65485.375
$ws.Range("I61").Value = 2903
$ws.Range("J61").Value = 169789.33
$ws.Range("K61").Value = 2903
$ws.Range("L61").Value = 169789.33
$ws.Range("M61").Value = -2691
$ws.Range("N61").Value = -170213.33

$ws.Range("H74").Value = 39493.555
$ws.Range("I74").Value = 57118.5
$ws.Range("J74").Value = 4243.6665
$ws.Range("K74").Value = 57118.5
$ws.Range("L74").Value = 4243.6665
$ws.Range("M74").Value = -56244.5
$ws.Range("N74").Value = -5991.6665

$ws.Range("H77").Value = 39493.555
$ws.Range("I77").Value = 57118.5
$ws.Range("J77").Value = 4243.6665
$ws.Range("K77").Value = 285592.5
$ws.Range("L77").Value = 21218.3325
$ws.Range("M77").Value = -281224.5
$ws.Range("N77").Value = -29954.3325

$ws.Range("H80").Value = 80000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 80000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 80000
$ws.Range("N80").Value = -81996

$ws.Range("H83").Value = 80000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 80000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 240000
$ws.Range("N83").Value = -249984

$ws.Range("H122").Value = 1831.76
$ws.Range("I122").Value = 1167.7646
$ws.Range("J122").Value = 3242.75
$ws.Range("K122").Value = 3503.2938
$ws.Range("L122").Value = 9728.25
$ws.Range("M122").Value = -1053.2938
$ws.Range("N122").Value = -14628.25

$ws.Range("H132").Value = 1749.9756
$ws.Range("I132").Value = 1623.1765
$ws.Range("J132").Value = 2365.8572
$ws.Range("K132").Value = 4869.529500000001
$ws.Range("L132").Value = 7097.571599999999
$ws.Range("M132").Value = -2339.529500000001
$ws.Range("N132").Value = -12157.5716

$ws.Range("H136").Value = 65485.375
$ws.Range("I136").Value = 2903
$ws.Range("J136").Value = 169789.33
$ws.Range("K136").Value = 8709
$ws.Range("L136").Value = 509367.99
$ws.Range("M136").Value = -6159
$ws.Range("N136").Value = -514467.99

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3469.0476
$ws.Range("I86").Value = 4111.3076
$ws.Range("J86").Value = 2425.375
$ws.Range("K86").Value = 4111.3076
$ws.Range("L86").Value = 2425.375
$ws.Range("M86").Value = -2988.3076
$ws.Range("N86").Value = -4671.375

$ws.Range("H89").Value = 3469.0476
$ws.Range("I89").Value = 4111.3076
$ws.Range("J89").Value = 2425.375
$ws.Range("K89").Value = 20556.538
$ws.Range("L89").Value = 12126.875
$ws.Range("M89").Value = -14940.538
$ws.Range("N89").Value = -23358.875

$ws.Range("H94").Value = 5701.6
$ws.Range("I94").Value = 4000.8572
$ws.Range("J94").Value = 9670
$ws.Range("K94").Value = 4000.8572
$ws.Range("L94").Value = 9670
$ws.Range("M94").Value = -3549.8572
$ws.Range("N94").Value = -10572

$ws.Range("H107").Value = 7145162
$ws.Range("I107").Value = 9093055
$ws.Range("J107").Value = 2887.6667
$ws.Range("K107").Value = 9093055
$ws.Range("L107").Value = 2887.6667
$ws.Range("M107").Value = -9091135
$ws.Range("N107").Value = -6727.6667

$ws.Range("H134").Value = 5137.773
$ws.Range("I134").Value = 2898.9412
$ws.Range("J134").Value = 12749.8
$ws.Range("K134").Value = 8696.8236
$ws.Range("L134").Value = 38249.39999999999
$ws.Range("M134").Value = -6161.8236
$ws.Range("N134").Value = -43319.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 299
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 299
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 299
$ws.Range("N6").Value = -525

$ws.Range("H7").Value = 190.58621
$ws.Range("I7").Value = 70.40000000000001
$ws.Range("J7").Value = 253.8421
$ws.Range("K7").Value = 70.40000000000001
$ws.Range("L7").Value = 253.8421
$ws.Range("M7").Value = 42.59999999999999
$ws.Range("N7").Value = -479.8421

$ws.Range("H17").Value = 516.6667
$ws.Range("I17").Value = 516.6667
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 516.6667
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -342.6667

$ws.Range("H22").Value = 215.25
$ws.Range("I22").Value = 187
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 187
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 163
$ws.Range("N22").Value = -1000

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()

$ws.Range("H41").Value = 13615.385
$ws.Range("I41").Value = 10500

$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -2251

$ws.Range("H69").Value = 47499
$ws.Range("I69").Value = 47499
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 47499
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -46750

$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9000
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -5256

$ws.Range("H72").Value = 47499
$ws.Range("I72").Value = 47499
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 142497
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -138753

$ws.Range("H86").Value = 4470648
$ws.Range("I86").Value = 8932296
$ws.Range("J86").Value = 8999.75
$ws.Range("K86").Value = 8932296
$ws.Range("L86").Value = 8999.75
$ws.Range("M86").Value = -8931173
$ws.Range("N86").Value = -11245.75

$ws.Range("H89").Value = 4470648
$ws.Range("I89").Value = 8932296
$ws.Range("J89").Value = 8999.75
$ws.Range("K89").Value = 44661480
$ws.Range("L89").Value = 44998.75
$ws.Range("M89").Value = -44655864
$ws.Range("N89").Value = -56230.75

$ws.Range("H103").Value = 1725
$ws.Range("I103").Value = 1725
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1725
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -553

$ws.Range("H107").Value = 1473.5333
$ws.Range("I107").Value = 1231
$ws.Range("J107").Value = 1790.6923
$ws.Range("K107").Value = 1231
$ws.Range("L107").Value = 1790.6923
$ws.Range("M107").Value = 689
$ws.Range("N107").Value = -5630.6923

$ws.Range("H132").Value = 3898859.2
$ws.Range("I132").Value = 5683435
$ws.Range("J132").Value = 1859344.1
$ws.Range("K132").Value = 17050305
$ws.Range("L132").Value = 5578032.300000001
$ws.Range("M132").Value = -17047775
$ws.Range("N132").Value = -5583092.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 84791.25
$ws.Range("I14").Value = 84791.25
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 254373.75
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -254200.75

$ws.Range("H38").Value = 131.5625
$ws.Range("I38").Value = 53
$ws.Range("J38").Value = 149.6923
$ws.Range("K38").Value = 159
$ws.Range("L38").Value = 449.0769
$ws.Range("M38").Value = 188
$ws.Range("N38").Value = -1143.0769

$ws.Range("H110").Value = 6731.75
$ws.Range("I110").Value = 6213.5
$ws.Range("J110").Value = 7250
$ws.Range("K110").Value = 18640.5
$ws.Range("L110").Value = 21750
$ws.Range("M110").Value = -14550.5
$ws.Range("N110").Value = -29930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3534.5
$ws.Range("I97").Value = 1750
$ws.Range("J97").Value = 4129.3335
$ws.Range("K97").Value = 1750
$ws.Range("L97").Value = 4129.3335
$ws.Range("M97").Value = -1254
$ws.Range("N97").Value = -5121.3335

$ws.Range("H113").Value = 3133477
$ws.Range("I113").Value = 224002.4
$ws.Range("J113").Value = 5558039
$ws.Range("K113").Value = 224002.4
$ws.Range("L113").Value = 5558039
$ws.Range("M113").Value = -221832.4
$ws.Range("N113").Value = -5562379

$ws.Range("H132").Value = 3420.7878
$ws.Range("I132").Value = 2954.1304
$ws.Range("J132").Value = 4494.1
$ws.Range("K132").Value = 8862.3912
$ws.Range("L132").Value = 13482.3
$ws.Range("M132").Value = -6332.3912
$ws.Range("N132").Value = -18542.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 235544.89
$ws.Range("I68").Value = 302286.28
$ws.Range("J68").Value = 1950
$ws.Range("K68").Value = 302286.28
$ws.Range("L68").Value = 1950
$ws.Range("M68").Value = -301537.28
$ws.Range("N68").Value = -3448

$ws.Range("H71").Value = 235544.89
$ws.Range("I71").Value = 302286.28
$ws.Range("J71").Value = 1950
$ws.Range("K71").Value = 1511431.4
$ws.Range("L71").Value = 9750
$ws.Range("M71").Value = -1507687.4
$ws.Range("N71").Value = -17238

$ws.Range("H133").Value = 35000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 35000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 35000
$ws.Range("N133").Value = -40060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5274014.5
$ws.Range("I96").Value = 15196
$ws.Range("J96").Value = 17544592
$ws.Range("K96").Value = 15196
$ws.Range("L96").Value = 17544592
$ws.Range("M96").Value = -13823
$ws.Range("N96").Value = -17547338

$ws.Range("H107").Value = 3620.5366
$ws.Range("I107").Value = 3961.6667
$ws.Range("J107").Value = 2213.375
$ws.Range("K107").Value = 11885.0001
$ws.Range("L107").Value = 6640.125
$ws.Range("M107").Value = -9965.000100000001
$ws.Range("N107").Value = -10480.125

$ws.Range("H122").Value = 2011.3684
$ws.Range("I122").Value = 1356.1818
$ws.Range("J122").Value = 2912.25
$ws.Range("K122").Value = 4068.5454
$ws.Range("L122").Value = 8736.75
$ws.Range("M122").Value = -1618.5454
$ws.Range("N122").Value = -13636.75

$ws.Range("H126").Value = 95146.664
$ws.Range("I126").Value = 139907.5
$ws.Range("J126").Value = 5625
$ws.Range("K126").Value = 419722.5
$ws.Range("L126").Value = 16875
$ws.Range("M126").Value = -417252.5
$ws.Range("N126").Value = -21815

$ws.Range("H132").Value = 2733.5454
$ws.Range("I132").Value = 2477.375
$ws.Range("J132").Value = 3416.6667
$ws.Range("K132").Value = 7432.125
$ws.Range("L132").Value = 10250.0001
$ws.Range("M132").Value = -4902.125
$ws.Range("N132").Value = -15310.0001
